# Update cryptos list cell values to match the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.355.14"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.592.86"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.48"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.507"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.815.56"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.577.04"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.67"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.345.16"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0733"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +4.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "212.25"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.03"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.08"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.27"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.335.81"
$ws.Range("E34").Value = "  +4.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.603"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.72"
$ws.Range("E40").Value = "  +4.62%  "
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.990"
$ws.Range("E42").Value = "  -24.82%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.766"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.728.70"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.90"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.31"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.49"
$ws.Range("E49").Value = "  -3.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0504"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0980"
$ws.Range("E51").Value = "  -4.21%  "
